$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 7.79
$ws.Range("F2").Value = 13.39
$ws.Range("N2").Value = 85.83574689470727
$ws.Range("N3").Value = 85.83574689470727
